# Update employee absence data rows 2-11 with new values as per the diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @{ Row=2;  A=14405; B="Alícia Pires";          C="Engenharia";             D="Problemas pessoais"; E=5; F=45082; G=2709.93 },
    @{ Row=3;  A=82326; B="Breno Rocha";            C="P&D";                    D="Problemas pessoais"; E=6; F=45096; G=3458.02 },
    @{ Row=4;  A=34256; B="Luiza Pacheco";           C="Financeiro";             D="Problemas pessoais"; E=1; F=45089; G=3361.77 },
    @{ Row=5;  A=47345; B="Maria Liz Moura";         C="Atendimento ao Cliente"; D="Viagem de negocios"; E=8; F=45088; G=8832.309999999999 },
    @{ Row=6;  A=6113;  B="Dr. Cauã Gomes";          C="Operacoes";              D="Problemas pessoais"; E=4; F=45089; G=9343.120000000001 },
    @{ Row=7;  A=84876; B="Dra. Cecília da Rosa";    C="Financeiro";             D="Doenca";              E=1; F=45092; G=4041.16 },
    @{ Row=8;  A=27234; B="Melina Dias";             C="P&D";                    D="Problemas pessoais"; E=3; F=45094; G=6078.2 },
    @{ Row=9;  A=2131;  B="Maria Isis Cassiano";     C="P&D";                    D="Outros";              E=6; F=45104; G=2349.19 },
    @{ Row=10; A=31710; B="Lavínia Aparecida";       C="Vendas";                 D="Outros";              E=5; F=45091; G=7535.09 },
    @{ Row=11; A=4192;  B="Marcela Costela";         C="Atendimento ao Cliente"; D="Problemas pessoais"; E=1; F=45106; G=9869.6 }
)

foreach ($item in $data) {
    $r = $item.Row
    $ws.Cells.Item($r, 1).Value = $item.A
    $ws.Cells.Item($r, 2).Value = $item.B
    $ws.Cells.Item($r, 3).Value = $item.C
    $ws.Cells.Item($r, 4).Value = $item.D
    $ws.Cells.Item($r, 5).Value = $item.E
    $ws.Cells.Item($r, 6).Value = $item.F
    $ws.Cells.Item($r, 7).Value = $item.G
}

$wb.Save()
